$wb = $excel.ActiveWorkbook

# --- Overview sheet: mark the second file (df08d9e8...) as ready for handoff ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("D3").Value = "2016-03-25 03:00:21"

# --- zh-cn sheet: update status + handoff datetime for the same file row ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "2016-03-25 03:00:15"

# --- de-de sheet: update status + handoff datetime for the same file row ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "2016-03-25 03:00:21"
